$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Docs.opencv.org. 2021. OpenCV: Basic concepts of the homography explained with code. [online] Available at: <https://docs.opencv.org/master/d9/dab/tutorial_homography.html> [Accessed 15 March 2021]."
$ws.Range("B12").Value = "(OpenCV: Basic concepts of the homography explained with code, 2021)"
$ws.Range("C12").Value = "results"
$ws.Range("D12").Value = "Detecting Social Distancing by  bounding box"

$ws.Range("B12").Select()
